$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B/C/D values for existing rows 2-17 (annotation_id, start_pos, end_pos)
$updates = @(
    ,@(2, 0, 28, 2307)
    ,@(3, 1, 25, 2298)
    ,@(4, 2, 119, 391)
    ,@(5, 3, 25, 2175)
    ,@(6, 4, 25, 594)
    ,@(7, 4, 596, 784)
    ,@(8, 5, 22, 69)
    ,@(9, 5, 70, 1059)
    ,@(10, 6, 1060, 1725)
    ,@(11, 7, 46, 1542)
    ,@(12, 8, 21, 1430)
    ,@(13, 9, 226, 784)
    ,@(14, 10, 26, 51)
    ,@(15, 10, 740, 1007)
    ,@(16, 11, 15, 707)
    ,@(17, 12, 15, 44)
)
foreach ($u in $updates) {
    $r = $u[0]
    $ws.Cells.Item($r, 2).Value = $u[1]
    $ws.Cells.Item($r, 3).Value = $u[2]
    $ws.Cells.Item($r, 4).Value = $u[3]
}

# Add new row 18, copying the formatting of row 17 (A:D and F) then setting its values
$ws.Range("A17:D17").Copy()
$ws.Range("A18:D18").PasteSpecial(-4122)
$ws.Range("F17").Copy()
$ws.Range("F18").PasteSpecial(-4122)
$ws.Range("A18").Value = 16
$ws.Range("B18").Value = 12
$ws.Range("C18").Value = 517
$ws.Range("D18").Value = 852

# Update the saved selection / active cell to D18 to mirror the authored view state
$ws.Range("D18").Select()
$excel.CutCopyMode = 0
